# Config.xlsx edit — "Code reorganization and bug fixes"
#  * Replace the stale "New Dataset" row on the Datasets sheet with the new
#    "Diagnosis" sample dataset row.
#  * Switch the active/selected sheet from "Methods" to "Datasets".

$wb = $excel.ActiveWorkbook

$datasets = $wb.Worksheets.Item("Datasets")

# --- Update the last row (row 14) of the Datasets sheet -------------------
# Old: "New Dataset" / ..\..\..\..\..\documents\Book1.xlsx / " Nothing"
# New: "Diagnosis "  / .\datasets\diagnosis_only.xlsx / description
$datasets.Range("B14").Value = "Diagnosis "
$datasets.Range("C14").Value = ".\datasets\diagnosis_only.xlsx"
$datasets.Range("D14").Value = 'A sample to test diagnosis without having an "actual" diagnosis column available.'

# --- Make "Datasets" the active tab/sheet (was "Methods") -----------------
$datasets.Activate()
